$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("A2").Value = 111950243
$ws.Range("B2").Value = 90823
$ws.Range("E2").Value = 5966
$ws.Range("F2").Value = "Motaggsvamp"
$ws.Range("G2").Value = "Sarcodon squamosus"
$ws.Range("H2").Value = "(Schaeff.) Quél."
$ws.Range("Q2").Value = 465473
$ws.Range("R2").Value = 6875785

# Row 3 updates
$ws.Range("A3").Value = 111950173
$ws.Range("B3").Value = 90792
$ws.Range("E3").Value = 4361
$ws.Range("F3").Value = "Orange taggsvamp"
$ws.Range("G3").Value = "Hydnellum aurantiacum"
$ws.Range("H3").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("Q3").Value = 465440
$ws.Range("R3").Value = 6875680
